# Logic to write API data to excel
# Writes a tenant/folder/trigger data table (as pulled from the Orchestrator
# API) into Sheet1: a header row plus two data rows describing the
# "Hackweek23" tenant's triggers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- Header row -----
$ws.Range("A1").Value = "TenantName"
$ws.Range("B1").Value = "FolderName"
$ws.Range("C1").Value = "Trigger_Process_Name"
$ws.Range("D1").Value = "Trigger_Schedule_CRON"
$ws.Range("E1").Value = "Trigger_Runtime_type"

# ----- Row 2: tenant/folder only -----
$ws.Range("A2").Value = "Hackweek23"
$ws.Range("B2").Value = "Shared"

# ----- Row 3: full trigger record -----
$ws.Range("A3").Value = "Hackweek23"
$ws.Range("B3").Value = "Hackweek 23 Testing"
$ws.Range("C3").Value = "AZURE.VM.Conenct"
$ws.Range("D3").Value = "0 0 0 1/1 * ? *"
$ws.Range("E3").Value = "Unattended"

# Size the columns to fit the data that was just written.
$ws.Range("A1:E3").EntireColumn.AutoFit()

# Leave the sheet selection on the second row, as in the source workbook.
[void]$ws.Range("A2:XFD2").Select()
